$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Portugal")

# Duplicate the Portugal sheet to the end of the workbook and rename it
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Slovakia"

# Update the market name and user story reference for the new sheet
$new.Range("B2").Value = "Slovakia Market"
$new.Range("B4").Value = "NGC-2930/T3236/T3235"

# The copied rows kept a stale cached auto-height from Portugal's narrower
# column layout; re-fit them so they fall back to the sheet default height
$new.Rows.Item(3).AutoFit()
$new.Rows.Item(4).AutoFit()
$new.Rows.Item(5).AutoFit()

# Selection on the new sheet is a single cell rather than the B4:B5 block
$new.Range("B4").Select() | Out-Null

# Leave the old sheet with a full-sheet selection and make the new sheet active
$src.Cells.Select() | Out-Null
$new.Activate()
$new.Range("B4").Select() | Out-Null

Write-Host "Added Slovakia worksheet"
